# nb: separating pilot 1 and pilot 2
# Re-label the "Cohort 1" rows (2-5) from pilot "1.x" to "6.x" and backfill
# the spike_interface columns (I:L) for row 2 and the spike_interface
# columns (H:L) for row 8, plus fix the eib_MD subject id on row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Cohort 1, Subject 1.1 -> 6.1) - fill in spike_interface_vHPC/BLA/LH/MD
$ws.Range("B2").Value = 6.1
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 16

# Row 3 (Cohort 1, Subject 1.2 -> 6.2)
$ws.Range("B3").Value = 6.2

# Row 4 (Cohort 1, Subject 1.3 -> 6.3)
$ws.Range("B4").Value = 6.3

# Row 5 (Cohort 1, Subject 1.4 -> 6.4)
$ws.Range("B5").Value = 6.4

# Row 8 (Cohort 2, Subject 1.3) - fill in spike_interface_mPFC/vHPC/BLA/LH/MD
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 31
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 29
$ws.Range("L8").Value = 28

# Row 9 (Cohort 2, Subject 1.4) - correct spike_interface_mPFC value
$ws.Range("H9").Value = 15

# Update the sheet's active cell / selection to match the author's last position
$ws.Range("H11").Select()
